$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source data formatting)
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"

$ws.Range('D2').Value = '36.805.28'
$ws.Range('E2').Value = '  -0.96%  '
$ws.Range('D3').Value = '2.091.78'
$ws.Range('E3').Value = '  +1.93%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '245.45'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('E6').Value = '  -1.57%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '54.31'
$ws.Range('E8').Value = '  -5.24%  '
$ws.Range('D9').Value = '58.98'
$ws.Range('E9').Value = '  -1.85%  '
$ws.Range('D10').Value = '0.367'
$ws.Range('E10').Value = '  -3.76%  '
$ws.Range('D11').Value = '0.0763'
$ws.Range('E11').Value = '  -2.22%  '
$ws.Range('E12').Value = '  +0.94%  '
$ws.Range('D13').Value = '0.915'
$ws.Range('E13').Value = '  +2.92%  '
$ws.Range('D14').Value = '15.06'
$ws.Range('E14').Value = '  -6.55%  '
$ws.Range('D15').Value = '2.396.66'
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('D16').Value = '5.51'
$ws.Range('E16').Value = '  -3.62%  '
$ws.Range('D17').Value = '2.122.97'
$ws.Range('E17').Value = '  +3.48%  '
$ws.Range('D18').Value = '36.792.97'
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').Value = '17.14'
$ws.Range('E19').Value = '  -6.57%  '
$ws.Range('D20').Value = '72.72'
$ws.Range('E20').Value = '  -2.79%  '
$ws.Range('D21').Value = '0.0₃0882'
$ws.Range('E21').Value = '  -1.32%  '
$ws.Range('D22').Value = '5.47'
$ws.Range('E22').Value = '  +1.46%  '
$ws.Range('D23').Value = '238.74'
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('E25').Value = '  -3.77%  '
$ws.Range('D26').Value = '9.75'
$ws.Range('E26').Value = '  +2.05%  '
$ws.Range('D27').Value = '2.17'
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('D28').Value = '166.80'
$ws.Range('E28').Value = '  -1.69%  '
$ws.Range('D29').Value = '20.90'
$ws.Range('E29').Value = '  +4.04%  '
$ws.Range('E30').Value = '  -1.55%  '
$ws.Range('D31').Value = '5.25'
$ws.Range('E31').Value = '  +8.40%  '
$ws.Range('E32').Value = '  +3.00%  '
$ws.Range('D33').Value = '4.70'
$ws.Range('E33').Value = '  +4.52%  '
$ws.Range('E34').Value = '  -1.57%  '
$ws.Range('D35').Value = '2.44'
$ws.Range('E35').Value = '  +8.04%  '
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').Value = '1.84'
$ws.Range('E37').Value = '  +3.64%  '
$ws.Range('E38').Value = '  -6.61%  '
$ws.Range('D39').Value = '1.27'
$ws.Range('E39').Value = '  -4.78%  '
$ws.Range('E40').Value = '  +1.34%  '
$ws.Range('E41').Value = '  -1.02%  '
$ws.Range('D42').Value = '4.89'
$ws.Range('E42').Value = '  -7.91%  '
$ws.Range('D43').Value = '0.0956'
$ws.Range('E43').Value = '  -3.92%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '96.43'
$ws.Range('E44').Value = '  +0.37%  '
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').Value = '2.84'
$ws.Range('E45').Value = '  -9.16%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.395.19'
$ws.Range('E46').Value = '  +10.06%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '16.07'
$ws.Range('E47').Value = '  -6.95%  '
$ws.Range('D48').Value = '7.47'
$ws.Range('E48').Value = '  +9.46%  '
$ws.Range('D49').Value = '2.46'
$ws.Range('E49').Value = '  +0.77%  '
$ws.Range('D50').Value = '2.91'
$ws.Range('E50').Value = '  +2.13%  '
$ws.Range('D51').Value = '2.284.85'
$ws.Range('E51').Value = '  +1.98%  '
